$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.703.43"
$ws.Range("E2").Value = "  +3.08%  "
$ws.Range("D3").Value = "1.789.41"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'223.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'0.555"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'32.57"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.98%  "
$ws.Range("D9").Value = "'0.280"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +3.78%  "
$ws.Range("D11").Value = "'0.0935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").Value = "2.046.45"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  +10.05%  "
$ws.Range("D14").Value = "1.783.43"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").Value = "34.702.89"
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "'253.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'10.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.54%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'158.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "'16.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "'3.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "1.430.28"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("D36").Value = "'1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.64%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.630"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0189"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").Value = "'82.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0504"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "'5.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.53%  "
$ws.Range("D47").Value = "1.942.71"
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("D48").Value = "'103.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.37%  "
$ws.Range("D49").Value = "'0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "'11.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("D51").Value = "'49.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.83%  "
